$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 30, shifting existing rows 30-63 down to 31-64
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new entry
$ws.Range("R30").Value = "hear your feedback atm"
$ws.Range("S30").Value = "2024-09-05 14:21:08"
